$p = $ppt.ActivePresentation

# Title textboxes use spAutoFit + wrap="none"; re-assigning their text makes
# the host recompute a tight bounding height immediately (real PowerPoint
# only does this on a render pass), which would otherwise shrink the shape.
# Restore the original authored height (57.6pt / 731520 EMU) after each
# title edit so the shape geometry is left untouched, matching the target.
$titleHeight = 57.60001

# ------------------------------------------------------------------
# Slide 2: "KEY FINANCIAL HIGHLIGHTS" -> "BUSINESS PROFILE & INFRASTRUCTURE"
# ------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Title
$t2 = $s2.Shapes.Item(2)
$t2.TextFrame.TextRange.Runs(1).Text = "BUSINESS PROFILE & INFRASTRUCTURE"
$t2.Height = $titleHeight

# Bullet body (shape 5) - drop last two bullets, rewrite remaining four
$tr2 = $s2.Shapes.Item(5).TextFrame.TextRange
$tr2.Paragraphs(2).Runs(1).Text = "■ The business operates in the entertainment sector."
$tr2.Paragraphs(3).Runs(1).Text = "■ It provides cinematic experiences across multiple locations."
$tr2.Paragraphs(4).Runs(1).Text = "■ The business has a significant physical presence in the region it operates."
$tr2.Paragraphs(5).Runs(1).Text = "■ It adheres to industry standards and certifications for operations."
# remove the two trailing bullets (ROE, Asset Turnover) - select from the end
# of paragraph 5 through (and past) the end of the text range so the whole
# tail, including its paragraph marks, is cleanly removed with no stray
# empty paragraph left behind.
$para5End = $tr2.Paragraphs(5).Start + $tr2.Paragraphs(5).Length
$tr2.Characters($para5End, $tr2.Length).Delete()

# ------------------------------------------------------------------
# Slide 3: "CASH FLOW ANALYSIS" -> "FINANCIAL & OPERATIONAL SCALE"
# ------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Title
$t3 = $s3.Shapes.Item(2)
$t3.TextFrame.TextRange.Runs(1).Text = "FINANCIAL & OPERATIONAL SCALE"
$t3.Height = $titleHeight

# Bullet body (shape 5) - rewrite first two bullets, append three more
$tr3 = $s3.Shapes.Item(5).TextFrame.TextRange
$tr3.Paragraphs(2).Runs(1).Text = "■ The business's Revenue From Operations increased significantly from 2022 to 2025."
$tr3.Paragraphs(3).Runs(1).Text = "■ Operating EBITDA also saw a significant increase, peaked in 2024 and slightly decreased in 2025."
[void]$tr3.InsertAfter("`r■ PAT Margin experienced fluctuations, with an increase observed in 2024 and 2025 after a decrease in the previous year.")
[void]$tr3.InsertAfter("`r■ ROE showed an upward trend, peaking in the same year as Operating EBITDA.")
[void]$tr3.InsertAfter("`r■ Asset Turnover ratio improved progressively over these years.")

# ------------------------------------------------------------------
# Slide 4: "KEY OBSERVATIONS" -> "INVESTMENT HIGHLIGHTS"
# ------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Title
$t4 = $s4.Shapes.Item(2)
$t4.TextFrame.TextRange.Runs(1).Text = "INVESTMENT HIGHLIGHTS"
$t4.Height = $titleHeight

# Bullet body (shape 5) - rewrite all four bullets in place
$tr4 = $s4.Shapes.Item(5).TextFrame.TextRange
$tr4.Paragraphs(2).Runs(1).Text = "■ The business has seen an increase in orders received over the years."
$tr4.Paragraphs(3).Runs(1).Text = "■ Capacity or production figures have been on a rise, indicating growth and expansion strategies."
$tr4.Paragraphs(4).Runs(1).Text = "■ The business maintains an order book with significant values across the years mentioned."
$tr4.Paragraphs(5).Runs(1).Text = "■ Management experience spans over several years, contributing to operational stability and growth."
